$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column H for rows 66 through 93 to the value 260105
for ($r = 66; $r -le 93; $r++) {
    $ws.Cells.Item($r, 8).Value = 260105
}

# Update the active selection from N86 to M86
$ws.Range("M86").Select()
